$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = [ordered]@{
    3  = "942546.03"
    4  = "1856615.58"
    5  = "103743.50"
    6  = "1870.33"
    7  = "6417.59"
    8  = "6875.05"
    9  = "19.15"
    10 = "971.99"
    11 = "12720.64"
    12 = "17.60"
    13 = "1669827.70"
    14 = "15779.13"
    15 = "137542.26"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("Q$row")
    $cell.Value = "'" + $values[$row]
}
